$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: new match data (08/08/2025 vs Estudiantes)
# Force the date column to stay plain text (dd/mm/yyyy, matching the rest
# of the sheet) instead of being auto-converted into a date serial number.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "08/08/2025"
$ws.Range("A21").Style = "Normal"
$ws.Range("B21").Value = "Estudiantes"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "Ind. Rivadavia"
$ws.Range("F21").Value = "L"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0.97
$ws.Range("L21").Value = 1.32
$ws.Range("M21").Value = 7
$ws.Range("N21").Value = 11
$ws.Range("O21").Value = 4
$ws.Range("P21").Value = 4
